$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Env_Staging")

# Update cells in the order that matches the new shared-string insertion order:
# 25, Male, khanna_deepankar@senger.co, DummyField, 11, Deepankar Khanna
$ws.Range("E3").Value = "25"
$ws.Range("C3").Value = "Male"
$ws.Range("B3").Value = "khanna_deepankar@senger.co"
$ws.Range("A7").Value = "DummyField"
$ws.Range("B7").Value = "DummyField"
$ws.Range("C7").Value = "11"
$ws.Range("A3").Value = "Deepankar Khanna"

# Update selection to match new active cell
$ws.Range("C16").Select()
